$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B144").Value = "Astrazeneca Ab"
$ws.Range("B145").Value = "Rooyan Darou"
$ws.Range("B152").Value = "Oncotec Pharma Produktion Gmbh"
$ws.Range("B153").Value = "Ipsen"
$ws.Range("B191").Value = "Hbm Pharma"
$ws.Range("B192").Value = "Takeda Italia S.p.a."
$ws.Range("B193").Value = "Sanofi Aventis"
$ws.Range("B194").Value = "Sydler"
$ws.Range("B195").Value = "Astrazeneca"
$ws.Range("B196").Value = "Excella Gmbh"
$ws.Range("B197").Value = "Farmigea"
$ws.Range("B198").Value = "Cipla Ltd"
$ws.Range("B199").Value = "Asa Daroo Toos"
$ws.Range("B200").Value = "Medochemie"
$ws.Range("B201").Value = "Merck Kgaa"
$ws.Range("B202").Value = "Desitin"
$ws.Range("B214").Value = "F.hoffmann-la Roche"
$ws.Range("B215").Value = "Imen Vaccine"
$ws.Range("B227").Value = "Rafarm"
$ws.Range("B228").Value = "Zydus Oncosciences"
$ws.Range("B229").Value = "Sun Pharmaceuticals Industries Ltd"
$ws.Range("B230").Value = "Generis"
$ws.Range("B231").Value = "Bausch & Lomb"
$ws.Range("B232").Value = "Nano Daru Pajuhan Pardis"
$ws.Range("B244").Value = "TURKTIPSAN A.S."
$ws.Range("B245").Value = "Johnlee Pharmaceuticals Private Limited"
$ws.Range("B246").Value = "Softgel Healthcare"
$ws.Range("B247").Value = "Bayer"
$ws.Range("B248").Value = "Hetero Labs Limited"
$ws.Range("B265").Value = "BioMarin Deutschland GmbH"
$ws.Range("B266").Value = "Karfarma"
$ws.Range("B284").Value = "Serb"
$ws.Range("B285").Value = "Tillomed Pharma GmbH"
$ws.Range("B286").Value = "polpharma sa"
$ws.Range("B287").Value = "Tillotts Pharma"
$ws.Range("B288").Value = "Sanofi-aventis S.p.a."
$ws.Range("B289").Value = "Amino Ag"
$ws.Range("B317").Value = "Alpex Pharma Sa"
$ws.Range("B318").Value = "S.c. Sandoz Srl"
$ws.Range("B319").Value = "Sanofi S.p.A."
$ws.Range("B320").Value = "Fareva (ex Pierre Fabre Medicament)"
$ws.Range("B321").Value = "Mylan Pharmaceuticals Inc"
$ws.Range("B322").Value = "Plus Pharma"
$ws.Range("B323").Value = "Wyeth Lederle Italia S.p.a"
$ws.Range("B324").Value = "Sanofi"
$ws.Range("B325").Value = "Fareva Mirabel (ex Merck Sharp & Dohme)"
$ws.Range("B326").Value = "Novonordisk"
$ws.Range("B327").Value = "Om Pharma"
$ws.Range("B328").Value = "Korea United Pharma"
$ws.Range("B351").Value = "Renaudin"
$ws.Range("B352").Value = "Virchow Drugs Limited"
$ws.Range("B353").Value = "Bayer Turk Kimya"
$ws.Range("B354").Value = "mahdaru Toba Nazarabad"
$ws.Range("B355").Value = "Jahanalcoholteb Arak"
$ws.Range("B356").Value = "GAP pharmaceuticals"
$ws.Range("B357").Value = "Catalent U.k. Swindon Zydis Limited"
$ws.Range("B358").Value = "Hk Pharma Limited"
$ws.Range("B359").Value = "Msd Schering Plough Labo"
$ws.Range("B360").Value = "Laboratoires Macors"
$ws.Range("B361").Value = "Aprazer Healthcare Private Limited"
$ws.Range("B362").Value = "Scott-Edil Pharmacia Ltd"
$ws.Range("B363").Value = "Aspen Bad Oldesloe"
$ws.Range("B364").Value = "Abbott Gmbh & Co. Kg"
$ws.Range("B365").Value = "Biogaran, S.A.S"
$ws.Range("B366").Value = "توسن دارو"
$ws.Range("B414").Value = "Ravenbhel Biotech"
$ws.Range("B415").Value = "Alcon-couvreur"
$ws.Range("B416").Value = "Reliance Industries Ltd"
$ws.Range("B417").Value = "Janssen-cilag"
$ws.Range("B418").Value = "Venus Remedies"
$ws.Range("B419").Value = "Mylan Laboratories Ltd"
$ws.Range("B420").Value = "Eli Lilly And Co"
$ws.Range("B421").Value = "Std Pharmaceutical Products Ltd"
$ws.Range("B422").Value = "Zydus Cadila"
$ws.Range("B423").Value = "Alphapharm Pty Ltd"
$ws.Range("B424").Value = "Dales Pharmaceuticals Limited"
$ws.Range("B425").Value = "Genzyme"
$ws.Range("B426").Value = "Aurovitas"
$ws.Range("B427").Value = "Concord Biotech Limited"
$ws.Range("B428").Value = "Haupt Pharma Wolfratshausen"
$ws.Range("B429").Value = "Intas"
$ws.Range("B430").Value = "Gufic Biosciences Ltd"
$ws.Range("B431").Value = "Holopack Verpackungstechnik Gmbh"
$ws.Range("B432").Value = "taghtirkhorasan"
$ws.Range("B433").Value = "PersisGen Par"
$ws.Range("B434").Value = "Hanlim Pharm Co Ltd"
$ws.Range("B435").Value = "Lupin Limited"
$ws.Range("B436").Value = "Lusomedicamenta"
$ws.Range("B476").Value = "Ali Raif Ilac Sanayi (aris)"
$ws.Range("B477").Value = "LABORATOIRES COLUXIA"
$ws.Range("B502").Value = "Bag Health Care Gmbh"
$ws.Range("B503").Value = "Laboratoire Europhartech"
$ws.Range("B504").Value = "Roche Pharma AG"
$ws.Range("B505").Value = "James Alexander Co"
$ws.Range("B506").Value = "I.E. Ulagay Ilac Sanayii Turk A.S."
$ws.Range("B507").Value = "CIRON DRUGS & PHARMACEUTICALS PVT. LTD."
$ws.Range("B508").Value = "Nipro Pharma"
$ws.Range("B509").Value = "Shimi Onsor Maryam"
$ws.Range("B510").Value = "Excella GmbH & Co. KG"
$ws.Range("B511").Value = "Natco Pharma Ltd."
$ws.Range("B512").Value = "Fleet laboratories Ltd"
$ws.Range("B513").Value = "Anfarm Hellas S.A."
$ws.Range("B514").Value = "Krka, D. D., Novo Mesto"
$ws.Range("B515").Value = "Snow Pharmaceuticals, LLC"
$ws.Range("B516").Value = "ALTAN PHARMACEUTICALS SA."
$ws.Range("B517").Value = "SHRI HARI PHARMACEUTICALS"
$ws.Range("B518").Value = "Beltapharm Spa"
$ws.Range("B519").Value = "Esteve Quimica S.a."
$ws.Range("B520").Value = "Gen Ilac"
$ws.Range("B521").Value = "Haupt Pharma Wulfing Gmbh"
$ws.Range("B522").Value = "Isu Abxis"
$ws.Range("B523").Value = "Metta Life Sciences Private Limited"
$ws.Range("B524").Value = "Hope Pharma"
$ws.Range("B525").Value = "Grifols Usa, Llc"
$ws.Range("B526").Value = "Dr. Franz Khler Chemie GmbH"
